# Seq_Diag_Dependency_2 — insert a new lifeline column ("Car") into the
# sequence diagram. The lifeline marks living in columns A/B/C that used to
# belong to other actors shift right (A->E, B->C, C->B) on the rows that
# belong to this diagram block (rows 2-21), freeing up column(s) for the
# newly-inserted lifeline.
#
# Excel has no native "move cell" primitive, so each logical move is done as
# copy-value + copy-format (re-using the existing style index instead of
# synthesizing a new one) followed by clearing the source cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlPasteFormats = -4122

function Move-Cell {
    # Relocates a cell's value+format from $src to $dst. The source cell is
    # left completely empty (no residual style), matching a cell that simply
    # no longer exists at that address.
    param([string]$src, [string]$dst)

    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial($xlPasteValues) | Out-Null
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Range($src).Clear() | Out-Null
}

function Move-CellKeepStyle {
    # Same relocation, but the source cell keeps its formatting (lifeline
    # border) and just becomes content-empty, since it's still part of the
    # diagram layout after the move.
    param([string]$src, [string]$dst)

    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial($xlPasteValues) | Out-Null
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Range($src).ClearContents() | Out-Null
}

# Rows 2-7: lifeline in column A shifts to column E.
Move-Cell "A2" "E2"
Move-Cell "A3" "E3"
Move-Cell "A4" "E4"
Move-Cell "A5" "E5"
Move-Cell "A6" "E6"
Move-Cell "A7" "E7"

# Rows 8-10: column B shifts to C, column A shifts to E.
Move-Cell "B8"  "C8"
Move-Cell "A8"  "E8"
Move-Cell "B9"  "C9"
Move-Cell "A9"  "E9"
Move-Cell "B10" "C10"
Move-Cell "A10" "E10"

# Row 11: only column A shifts to E.
Move-Cell "A11" "E11"

# Rows 12-17: column C shifts to B, column A shifts to E.
Move-Cell "C12" "B12"
Move-Cell "A12" "E12"
Move-Cell "C13" "B13"
Move-Cell "A13" "E13"
Move-Cell "C14" "B14"
Move-Cell "A14" "E14"
Move-Cell "C15" "B15"
Move-Cell "A15" "E15"
Move-Cell "C16" "B16"
Move-Cell "A16" "E16"
Move-Cell "C17" "B17"
Move-Cell "A17" "E17"

# Row 18: column C shifts to B; the activation mark that used to sit in E
# moves back into A (E remains part of the layout, now content-empty).
Move-Cell "C18" "B18"
Move-CellKeepStyle "E18" "A18"

# Rows 19-21: column C shifts to B.
Move-Cell "C19" "B19"
Move-Cell "C20" "B20"
Move-Cell "C21" "B21"
